# Auto-generated from the unified OOXML diff.
# Updates the cryptos worksheet's Price (D) and Volume(1h) (E) columns,
# including two pairs of rows whose data swapped order (12/13, 39/40, 48/49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold text that often looks numeric (e.g. "26.302.11",
# "0.5170") or percentage-like strings with padding spaces
# (e.g. "  -2.96%  "). The source workbook stores all of these as plain
# text (inlineStr), so force the range to Text format first; otherwise
# Excel's COM layer will auto-coerce values like "0.5170" into the number
# 0.517 and drop the information the diff expects us to preserve.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.302.11'
$ws.Range("E2").Value = '  -2.96%  '
$ws.Range("D3").Value = '1.831.14'
$ws.Range("E3").Value = '  -2.61%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '258.65'
$ws.Range("E5").Value = '  -7.93%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").Value = '0.5170'
$ws.Range("E7").Value = '  -2.35%  '
$ws.Range("E8").Value = '  -8.77%  '
$ws.Range("D9").Value = '0.06720'
$ws.Range("E9").Value = '  -4.60%  '
$ws.Range("D10").Value = '18.61'
$ws.Range("E10").Value = '  -8.50%  '
$ws.Range("D11").Value = '0.7623'
$ws.Range("E11").Value = '  -6.73%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.888.49'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07693'
$ws.Range("E13").Value = '  -1.52%  '
$ws.Range("D14").Value = '88.44'
$ws.Range("E14").Value = '  -2.29%  '
$ws.Range("D15").Value = '5.007'
$ws.Range("E15").Value = '  -3.70%  '
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("E17").Value = '  -3.78%  '
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("E19").Value = '  -3.87%  '
$ws.Range("D20").Value = '26.320.83'
$ws.Range("E20").Value = '  -2.99%  '
$ws.Range("D21").Value = '2.077.56'
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("D22").Value = '4.524'
$ws.Range("E22").Value = '  -5.02%  '
$ws.Range("D23").Value = '9.395'
$ws.Range("E23").Value = '  -7.28%  '
$ws.Range("D24").Value = '5.889'
$ws.Range("E24").Value = '  -5.24%  '
$ws.Range("D25").Value = '2.300'
$ws.Range("E25").Value = '  -3.76%  '
$ws.Range("D26").Value = '144.58'
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("D27").Value = '1.641'
$ws.Range("E27").Value = '  -1.82%  '
$ws.Range("D28").Value = '16.88'
$ws.Range("E28").Value = '  -3.93%  '
$ws.Range("D29").Value = '110.61'
$ws.Range("E29").Value = '  -1.70%  '
$ws.Range("D30").Value = '4.170'
$ws.Range("E30").Value = '  -5.16%  '
$ws.Range("D31").Value = '4.105'
$ws.Range("E31").Value = '  -6.30%  '
$ws.Range("D32").Value = '0.08717'
$ws.Range("E32").Value = '  -2.28%  '
$ws.Range("D33").Value = '0.04816'
$ws.Range("E33").Value = '  -1.77%  '
$ws.Range("D34").Value = '1.122'
$ws.Range("E34").Value = '  -4.66%  '
$ws.Range("D35").Value = '2.854'
$ws.Range("E35").Value = '  -1.55%  '
$ws.Range("D36").Value = '0.6768'
$ws.Range("E36").Value = '  -9.39%  '
$ws.Range("D37").Value = '3.087'
$ws.Range("E37").Value = '  -6.62%  '
$ws.Range("D38").Value = '0.01776'
$ws.Range("E38").Value = '  -5.61%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '0.4881'
$ws.Range("E39").Value = '  -8.34%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '2.175'
$ws.Range("E40").Value = '  -10.08%  '
$ws.Range("D41").Value = '110.78'
$ws.Range("E41").Value = '  -5.37%  '
$ws.Range("D42").Value = '0.8894'
$ws.Range("E42").Value = '  -8.53%  '
$ws.Range("D43").Value = '6.109'
$ws.Range("E43").Value = '  -3.38%  '
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").Value = '7.660'
$ws.Range("E45").Value = '  -6.82%  '
$ws.Range("D46").Value = '0.4167'
$ws.Range("E46").Value = '  -9.36%  '
$ws.Range("E47").Value = '  -8.67%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05877'
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.995'
$ws.Range("E49").Value = '  -5.19%  '
$ws.Range("D50").Value = '35.18'
$ws.Range("E50").Value = '  -4.18%  '
$ws.Range("D51").Value = '59.03'
$ws.Range("E51").Value = '  -4.43%  '
